$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the trailing rows (25-27) that are removed in the target layout.
$ws.Range("A25:D27").Clear()

# Rewrite all data rows (1-24) with the reorganised/deduplicated dataset.
$ws.Range("A1").Value = 5272804613
$ws.Range("B1").Value = "0B50559"
$ws.Range("C1").Value = 201
$ws.Range("D1").Value = "ZW03"
$ws.Range("A2").Value = 5272804613
$ws.Range("B2").Value = "4Z51L42650"
$ws.Range("C2").Value = 201
$ws.Range("D2").Value = "ZC06"
$ws.Range("A3").Value = 5272813554
$ws.Range("B3").Value = "7D6MNCM3WW"
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = "ZC06"
$ws.Range("A4").Value = 5272813581
$ws.Range("B4").Value = "7D76CTO1WW"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = "ZC06"
$ws.Range("A5").Value = 5272815155
$ws.Range("B5").Value = "MOUSE00"
$ws.Range("C5").Value = 201
$ws.Range("D5").Value = "VIRTUAIS"
$ws.Range("A6").Value = 5272815155
$ws.Range("B6").Value = "TECLA00"
$ws.Range("C6").Value = 201
$ws.Range("D6").Value = "VIRTUAIS"
$ws.Range("A7").Value = 5272815155
$ws.Range("B7").Value = "13E0S00400"
$ws.Range("C7").Value = 21
$ws.Range("D7").Value = "ZC06"
$ws.Range("A8").Value = 5272815181
$ws.Range("B8").Value = "13E0S00400"
$ws.Range("C8").Value = 150
$ws.Range("D8").Value = "ZAGNA0902"
$ws.Range("A9").Value = 5272815181
$ws.Range("B9").Value = "13E0S00400"
$ws.Range("C9").Value = 150
$ws.Range("D9").Value = "ZAGNJ2204"
$ws.Range("A10").Value = 5272815181
$ws.Range("B10").Value = "TECLA00"
$ws.Range("C10").Value = 150
$ws.Range("D10").Value = "VIRTUAIS"
$ws.Range("A11").Value = 5272815181
$ws.Range("B11").Value = "MOUSE00"
$ws.Range("C11").Value = 150
$ws.Range("D11").Value = "VIRTUAIS"
$ws.Range("A12").Value = 5272815187
$ws.Range("B12").Value = "0B50559"
$ws.Range("C12").Value = 150
$ws.Range("D12").Value = "ZW07"
$ws.Range("A13").Value = 5272815187
$ws.Range("B13").Value = "4Z51L42650"
$ws.Range("C13").Value = 150
$ws.Range("D13").Value = "ZW07"
$ws.Range("A14").Value = 5272815335
$ws.Range("B14").Value = "63CFMAR1BR"
$ws.Range("C14").Value = 30
$ws.Range("D14").Value = "ZAGNW3203"
$ws.Range("A15").Value = 5272815360
$ws.Range("B15").Value = "21NQA07BBR"
$ws.Range("C15").Value = 125
$ws.Range("D15").Value = "ZAGNM0103"
$ws.Range("A16").Value = 5272815360
$ws.Range("B16").Value = "21NQA07BBR"
$ws.Range("C16").Value = 125
$ws.Range("D16").Value = "ZC06"
$ws.Range("A17").Value = 5272815362
$ws.Range("B17").Value = "21NQA07BBR"
$ws.Range("C17").Value = 126
$ws.Range("D17").Value = "ZAGNG2306"
$ws.Range("A18").Value = 5272815362
$ws.Range("B18").Value = "21NQA07BBR"
$ws.Range("C18").Value = 126
$ws.Range("D18").Value = "ZC06"
$ws.Range("A19").Value = 5272815364
$ws.Range("B19").Value = "0B50559"
$ws.Range("C19").Value = 150
$ws.Range("D19").Value = "ZW04"
$ws.Range("A20").Value = 5272815364
$ws.Range("B20").Value = "4Z51L42650"
$ws.Range("C20").Value = 150
$ws.Range("D20").Value = "ZW04"
$ws.Range("A21").Value = 5272815384
$ws.Range("B21").Value = "MOUSE00"
$ws.Range("C21").Value = 200
$ws.Range("D21").Value = "VIRTUAIS"
$ws.Range("A22").Value = 5272815384
$ws.Range("B22").Value = "TECLA00"
$ws.Range("C22").Value = 200
$ws.Range("D22").Value = "VIRTUAIS"
$ws.Range("A23").Value = 5272815655
$ws.Range("B23").Value = "63CFMAR1BR"
$ws.Range("C23").Value = 30
$ws.Range("D23").Value = "ZAGNW3803"
$ws.Range("A24").Value = 5272816410
$ws.Range("B24").Value = "21NQA0K6BR"
$ws.Range("C24").Value = 6
$ws.Range("D24").Value = "ZW06"

# Restore the reported active selection.
[void]$ws.Range("D8").Select()
